# Odds-base update (19-06-2024 21:51): for a number of fixtures the two
# rows (adjacent match rows) had their data swapped - i.e. the "id" and all
# match statistics/odds columns (B:AD) of each pair of rows were exchanged,
# while the leading rank column (A) stayed attached to its own row.
#
# Row pairs whose B:AD contents must be swapped:
$pairs = @(
    @(76, 77),
    @(84, 85),
    @(108, 109),
    @(117, 118),
    @(130, 131),
    @(133, 134),
    @(150, 151),
    @(159, 160),
    @(164, 165),
    @(170, 171),
    @(221, 222),
    @(226, 227)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $range1 = $ws.Range("B$row1`:AD$row1")
    $range2 = $ws.Range("B$row2`:AD$row2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
